# Update gh-pages to output generated at 456a3b4
#
# Both "展览" and "全部类型" sheets dropped the "南宁·国乙only" (2024-08-10)
# listing. Every remaining row shifts up by one, the running index in
# column A is renumbered, and three listings' "想去人数" (F column) counts
# were refreshed by the scraper:
#   熊喵M动漫嘉年华   1278 -> 1281
#   第二届北极光动漫展 2779 -> 2787
#   万圣漫控嘉年华10   247  -> 248

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the obsolete "南宁·国乙only" row (row 2); everything below
    # shifts up automatically.
    $ws.Range("A2").EntireRow.Delete()

    # Renumber the running index in column A (0 = header, 1..n = data).
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh "想去人数" (column F) for the three listings whose counts
    # the scraper picked up as changed.
    for ($r = 2; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value
        if ($name -eq "南宁·熊喵M动漫嘉年华【免费】") {
            $ws.Cells.Item($r, 6).Value = 1281
        } elseif ($name -eq "南宁·第二届北极光动漫展") {
            $ws.Cells.Item($r, 6).Value = 2787
        } elseif ($name -eq "南宁·万圣漫控嘉年华10") {
            $ws.Cells.Item($r, 6).Value = 248
        }
    }
}
